# Update market/profit figures across the Sheets workbook (scheduled runner sync).
# Generated from the authoritative diff of Asura_Profits sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 14059.6
$ws.Range("I8").Value = 99.333336
$ws.Range("J8").Value = 35000
$ws.Range("K8").Value = 298.000008
$ws.Range("L8").Value = 105000
$ws.Range("M8").Value = -159.000008
$ws.Range("N8").Value = -105278

$ws.Range("H31").Value = 1699316.5
$ws.Range("I31").Value = 1699316.5
$ws.Range("K31").Value = 5097949.5
$ws.Range("M31").Value = -5097719.5

$ws.Range("H76").Value = 3789.5
$ws.Range("I76").Value = 3643.5715
$ws.Range("J76").Value = 4130
$ws.Range("K76").Value = 3643.5715
$ws.Range("L76").Value = 4130
$ws.Range("M76").Value = -3328.5715
$ws.Range("N76").Value = -4760

$ws.Range("H79").Value = 3789.5
$ws.Range("I79").Value = 3643.5715
$ws.Range("J79").Value = 4130
$ws.Range("K79").Value = 3643.5715
$ws.Range("L79").Value = 4130
$ws.Range("M79").Value = -2551.5715
$ws.Range("N79").Value = -6314

$ws.Range("H94").Value = 5210.0625
$ws.Range("I94").Value = 5370.7334
$ws.Range("K94").Value = 5370.7334
$ws.Range("M94").Value = -4919.7334

$ws.Range("H100").Value = 2328
$ws.Range("I100").Value = 1302
$ws.Range("J100").Value = 2841
$ws.Range("K100").Value = 1302
$ws.Range("L100").Value = 2841
$ws.Range("M100").Value = -761
$ws.Range("N100").Value = -3923

$ws.Range("H132").Value = 2233.5833
$ws.Range("I132").Value = 1363.9354
$ws.Range("J132").Value = 3819.4119
$ws.Range("K132").Value = 4091.8062
$ws.Range("L132").Value = 11458.2357
$ws.Range("M132").Value = -1561.8062
$ws.Range("N132").Value = -16518.2357

$ws.Range("H138").Value = 3128169.5
$ws.Range("J138").Value = 4121.9756
$ws.Range("L138").Value = 12365.9268
$ws.Range("N138").Value = -22645.9268

$ws.Range("H141").Value = 11494.032
$ws.Range("I141").Value = 6009.3184
$ws.Range("J141").Value = 24901.111
$ws.Range("K141").Value = 18027.9552
$ws.Range("L141").Value = 74703.333
$ws.Range("M141").Value = -12847.9552
$ws.Range("N141").Value = -85063.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518
$ws.Range("M23").ClearContents()

$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496

$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716

$ws.Range("H102").Value = 1702
$ws.Range("I102").Value = 1836.6666
$ws.Range("K102").Value = 1836.6666
$ws.Range("M102").Value = -214.6666

$ws.Range("H109").Value = 32000
$ws.Range("J109").Value = 32000
$ws.Range("L109").Value = 32000
$ws.Range("N109").Value = -34774

$ws.Range("H110").Value = 1282
$ws.Range("I110").Value = 1282
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1282
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 763
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 13657.143
$ws.Range("I8").Value = 13657.143
$ws.Range("K8").Value = 13657.143
$ws.Range("M8").Value = -13517.143

$ws.Range("H62").Value = 29900
$ws.Range("I62").Value = 29800
$ws.Range("K62").Value = 29800
$ws.Range("M62").Value = -29114

$ws.Range("H65").Value = 29900
$ws.Range("I65").Value = 29800
$ws.Range("K65").Value = 89400
$ws.Range("M65").Value = -85968

$ws.Range("H86").Value = 101206.1
$ws.Range("I86").Value = 1391.6666
$ws.Range("J86").Value = 250927.75
$ws.Range("K86").Value = 1391.6666
$ws.Range("L86").Value = 250927.75
$ws.Range("M86").Value = -268.6666
$ws.Range("N86").Value = -253173.75

$ws.Range("H89").Value = 101206.1
$ws.Range("I89").Value = 1391.6666
$ws.Range("J89").Value = 250927.75
$ws.Range("K89").Value = 6958.333000000001
$ws.Range("L89").Value = 1254638.75
$ws.Range("M89").Value = -1342.333000000001
$ws.Range("N89").Value = -1265870.75

$ws.Range("H94").Value = 2350
$ws.Range("I94").Value = 1700
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1700
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -1249
$ws.Range("N94").Value = -3902

$ws.Range("H105").Value = 3379.2354
$ws.Range("I105").Value = 3238.077
$ws.Range("K105").Value = 3238.077
$ws.Range("M105").Value = -1491.077

$ws.Range("H108").Value = 58614
$ws.Range("J108").Value = 58614
$ws.Range("L108").Value = 58614
$ws.Range("N108").Value = -66294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1844.7142
$ws.Range("I94").Value = 973.8
$ws.Range("J94").Value = 2116.875
$ws.Range("K94").Value = 973.8
$ws.Range("L94").Value = 2116.875
$ws.Range("M94").Value = -522.8
$ws.Range("N94").Value = -3018.875

$ws.Range("H132").Value = 2023.125
$ws.Range("I132").Value = 1598
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4794
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2264
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 282.4
$ws.Range("I17").Value = 256
$ws.Range("K17").Value = 768
$ws.Range("M17").Value = -599

$ws.Range("H46").Value = 2119.6
$ws.Range("I46").Value = 199.66667
$ws.Range("J46").Value = 4999.5
$ws.Range("K46").Value = 599.00001
$ws.Range("L46").Value = 14998.5
$ws.Range("M46").Value = -508.00001
$ws.Range("N46").Value = -15180.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 37500100
$ws.Range("I7").Value = 37500100
$ws.Range("K7").Value = 37500100
$ws.Range("M7").Value = -37499988

$ws.Range("H8").Value = 37500100
$ws.Range("I8").Value = 37500100
$ws.Range("K8").Value = 37500100
$ws.Range("M8").Value = -37499961

$ws.Range("H80").Value = 3124.25
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3198.8
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3198.8
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5194.8

$ws.Range("H83").Value = 3124.25
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3198.8
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15994
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -25978

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1700
$ws.Range("I46").Value = 1750
$ws.Range("J46").Value = 1666.6666
$ws.Range("K46").Value = 1750
$ws.Range("L46").Value = 1666.6666
$ws.Range("M46").Value = -1562
$ws.Range("N46").Value = -2042.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 17101
$ws.Range("I8").Value = 751.5
$ws.Range("J8").Value = 49800
$ws.Range("K8").Value = 751.5
$ws.Range("L8").Value = 49800
$ws.Range("M8").Value = -611.5
$ws.Range("N8").Value = -50080

$ws.Range("H62").Value = 4207.6924
$ws.Range("I62").Value = 3650
$ws.Range("J62").Value = 4685.7144
$ws.Range("K62").Value = 3650
$ws.Range("L62").Value = 4685.7144
$ws.Range("M62").Value = -3026
$ws.Range("N62").Value = -5933.7144

$ws.Range("H65").Value = 4207.6924
$ws.Range("I65").Value = 3650
$ws.Range("J65").Value = 4685.7144
$ws.Range("K65").Value = 18250
$ws.Range("L65").Value = 23428.572
$ws.Range("M65").Value = -15130
$ws.Range("N65").Value = -29668.572

